# Insert a new row at position 51 (pushing existing rows 51..120 down to 52..121)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the new record's data
$ws.Range("A51").Value = 11
$ws.Range("B51").Value = "Vega Monumental Concepción"
$ws.Range("C51").Value = "Bíobío"
$ws.Range("D51").Value = 44771
$ws.Range("E51").Value = 8
$ws.Range("F51").Value = 100112021
$ws.Range("G51").Value = "Ají"
$ws.Range("H51").Value = "Inferno"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 100
$ws.Range("K51").Value = 15000
$ws.Range("L51").Value = 16000
$ws.Range("M51").Value = 15500
$ws.Range("N51").Value = "`$/caja 12 kilos"
$ws.Range("O51").Value = "Región de Arica y Parinacota"
$ws.Range("P51").Value = 1292
$ws.Range("Q51").Value = 12
$ws.Range("R51").Value = "Hortaliza"
